$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").ClearContents()

$ws.Range("C2").Value = 0.59215529571917125
$ws.Range("D2").Value = 0.18082567314225537
$ws.Range("E2").Value = 1.0600421034538978

$ws.Range("B3").Value = 0.13232252769231845
$ws.Range("C3").Value = 1.8279932993459795
$ws.Range("D3").Value = 0.73773049294567783
$ws.Range("E3").Value = 2.3739223718691895

$ws.Range("B1:E3").Select()
